# Apply vocabulary updates described in the commit's diff.
# The sheet is a flat SKOS-style vocabulary table with columns:
#   A: Identifier   B: prefLabel   C: altLabel   D: definition
#   E: source       F: broader     G: exactMatch
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Several child concepts of "animal" (id-amr:10026) are repointed to the
#    new generic "animal" concept gen:10025.
$animalChildren = @(68, 82, 83, 84, 86)
foreach ($r in $animalChildren) {
    $ws.Range("F$r").Value = "gen:10025"
}

# 2) Rows 105-110 are replaced with a new set of rows coming from the
#    regenerated Google Sheet / .ttl source. The old "animal" (id-amr:10026)
#    row is dropped, and the following rows shift up by one, with a brand
#    new "parasite" row appended at the end (row 110).

# Row 105: covid19:10151 / bacteria
$ws.Range("A105").Value = "covid19:10151"
$ws.Range("B105").Value = "bacteria"
$ws.Range("C105").Value = ""
$ws.Range("D105").Value = "Unicellular, prokaryotic organisms that reproduce by cell division and usually have cell walls; can be shaped like spheres, rods or spirals and can be found in virtually any environment."
$ws.Range("E105").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14187"
$ws.Range("F105").Value = "gen:10024"
$ws.Range("G105").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14187"

# Row 106: covid19:10150 / virus
$ws.Range("A106").Value = "covid19:10150"
$ws.Range("B106").Value = "virus"
$ws.Range("D106").Value = "An infectious agent which consists of two parts, genetic material and a protein coat. These organisms lack independent metabolism, and they must infect the cells of other types of organisms to reproduce. Most viruses are capable of passing through fine filters that retain bacteria, and are not visible through a light microscope."
$ws.Range("E106").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14283"
$ws.Range("F106").Value = "gen:10024"
$ws.Range("G106").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C14283"

# Row 107: gen:10011 / macromolecule
$ws.Range("A107").Value = "gen:10011"
$ws.Range("B107").Value = "macromolecule"
$ws.Range("D107").Value = ""
$ws.Range("E107").Value = ""
$ws.Range("F107").Value = "gen:10007"
$ws.Range("G107").Value = ""

# Row 108: id-amr:10085 / protein
$ws.Range("A108").Value = "id-amr:10085"
$ws.Range("B108").Value = "protein"
$ws.Range("F108").Value = "gen:10011"

# Row 109: id-amr:10086 / antibody
$ws.Range("A109").Value = "id-amr:10086"
$ws.Range("B109").Value = "antibody"
$ws.Range("F109").Value = "gen:10011"

# Row 110: covid19:10152 / parasite (brand new row)
$ws.Range("A110").Value = "covid19:10152"
$ws.Range("B110").Value = "parasite"
$ws.Range("D110").Value = "Any organism that has a close, symbiotic relationship with a separate, host organism."
$ws.Range("E110").Value = "https://ncithesaurus.nci.nih.gov/ncitbrowser/ConceptReport.jsp?dictionary=NCI_Thesaurus&ns=ncit&code=C28176"
$ws.Range("F110").Value = "gen:10024"
$ws.Range("G110").Value = "http://purl.obolibrary.org/obo/NCIT_C28176"
